# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Goblin_Profits workbook, per sheet (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 142857280
$ws.Range("J2").Value = 500000060
$ws.Range("L2").Value = 500000060
$ws.Range("N2").Value = -500000286
$ws.Range("H18").Value = 2102.7334
$ws.Range("I18").Value = 2221.2856
$ws.Range("J18").Value = 1999
$ws.Range("K18").Value = 2221.2856
$ws.Range("L18").Value = 1999
$ws.Range("M18").Value = -1937.2856
$ws.Range("N18").Value = -2567
$ws.Range("H98").Value = 7084.6665
$ws.Range("I98").Value = 6744.7896
$ws.Range("J98").Value = 20000
$ws.Range("K98").Value = 6744.7896
$ws.Range("L98").Value = 20000
$ws.Range("M98").Value = -5246.7896
$ws.Range("N98").Value = -22996
$ws.Range("H106").Value = 6595.6
$ws.Range("I106").Value = 6826.3335
$ws.Range("K106").Value = 6826.3335
$ws.Range("M106").Value = -6195.3335
$ws.Range("H112").Value = 2718
$ws.Range("J112").Value = 3205.2
$ws.Range("L112").Value = 9615.599999999999
$ws.Range("N112").Value = -11831.6
$ws.Range("H116").Value = 8983.846
$ws.Range("J116").Value = 8724
$ws.Range("L116").Value = 8724
$ws.Range("N116").Value = -15608
$ws.Range("H122").Value = 7084.6665
$ws.Range("I122").Value = 6744.7896
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 20234.3688
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -17784.3688
$ws.Range("N122").Value = -64900
$ws.Range("H132").Value = 3706471.2
$ws.Range("I132").Value = 2871.0908
$ws.Range("J132").Value = 20002312
$ws.Range("K132").Value = 8613.2724
$ws.Range("L132").Value = 60006936
$ws.Range("M132").Value = -6083.2724
$ws.Range("N132").Value = -60011996

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 68208.94500000001
$ws.Range("I32").Value = 69115.75
$ws.Range("K32").Value = 69115.75
$ws.Range("M32").Value = -68828.75
$ws.Range("H45").Value = 3164.9473
$ws.Range("I45").Value = 3175.9333
$ws.Range("K45").Value = 3175.9333
$ws.Range("M45").Value = -2798.9333
$ws.Range("H122").Value = 4631264.5
$ws.Range("I122").Value = 6174247
$ws.Range("J122").Value = 2316.5
$ws.Range("K122").Value = 18522741
$ws.Range("L122").Value = 6949.5
$ws.Range("M122").Value = -18520291
$ws.Range("N122").Value = -11849.5
$ws.Range("H139").Value = 73833.336
$ws.Range("J139").Value = 73833.336
$ws.Range("L139").Value = 73833.336
$ws.Range("N139").Value = -84113.336

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2158.56
$ws.Range("I20").Value = 2132.85
$ws.Range("J20").Value = 2261.4
$ws.Range("K20").Value = 2132.85
$ws.Range("L20").Value = 2261.4
$ws.Range("M20").Value = -1885.85
$ws.Range("N20").Value = -2755.4
$ws.Range("H94").Value = 3227.2307
$ws.Range("I94").Value = 2595.4
$ws.Range("J94").Value = 5333.3335
$ws.Range("K94").Value = 2595.4
$ws.Range("L94").Value = 5333.3335
$ws.Range("M94").Value = -2144.4
$ws.Range("N94").Value = -6235.3335
$ws.Range("H105").Value = 3115.3333
$ws.Range("I105").Value = 1954.25
$ws.Range("K105").Value = 1954.25
$ws.Range("M105").Value = -207.25
$ws.Range("H107").Value = 6211.4614
$ws.Range("I107").Value = 4035.5715
$ws.Range("J107").Value = 8750
$ws.Range("K107").Value = 4035.5715
$ws.Range("L107").Value = 8750
$ws.Range("M107").Value = -2115.5715
$ws.Range("N107").Value = -12590
$ws.Range("H134").Value = 1318873.8
$ws.Range("I134").Value = 2225.0625
$ws.Range("J134").Value = 8341000
$ws.Range("K134").Value = 6675.1875
$ws.Range("L134").Value = 25023000
$ws.Range("M134").Value = -4140.1875
$ws.Range("N134").Value = -25028070

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1899.75
$ws.Range("I99").Value = 1931.7
$ws.Range("J99").Value = 1740
$ws.Range("K99").Value = 1931.7
$ws.Range("L99").Value = 1740
$ws.Range("M99").Value = -433.7
$ws.Range("N99").Value = -4736
$ws.Range("H105").Value = 7057.6665
$ws.Range("I105").Value = 5074.143
$ws.Range("K105").Value = 5074.143
$ws.Range("M105").Value = -3327.143
$ws.Range("H126").Value = 1899.75
$ws.Range("I126").Value = 1931.7
$ws.Range("J126").Value = 1740
$ws.Range("K126").Value = 5795.1
$ws.Range("L126").Value = 5220
$ws.Range("M126").Value = -3325.1
$ws.Range("N126").Value = -10160
$ws.Range("H132").Value = 3201.8
$ws.Range("I132").Value = 2032.909
$ws.Range("J132").Value = 6416.25
$ws.Range("K132").Value = 6098.727000000001
$ws.Range("L132").Value = 19248.75
$ws.Range("M132").Value = -3568.727000000001
$ws.Range("N132").Value = -24308.75
$ws.Range("H141").Value = 189935.08
$ws.Range("J141").Value = 189935.08
$ws.Range("L141").Value = 189935.08
$ws.Range("N141").Value = -200295.08

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3091447.5
$ws.Range("I11").Value = 4956450.5
$ws.Range("J11").Value = 39624.727
$ws.Range("K11").Value = 4956450.5
$ws.Range("L11").Value = 39624.727
$ws.Range("M11").Value = -4956311.5
$ws.Range("N11").Value = -39902.727
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H45").Value = 74999
$ws.Range("J45").Value = 74999
$ws.Range("L45").Value = 74999
$ws.Range("N45").Value = -76117
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H80").Value = 83338750
$ws.Range("I80").Value = 200002600
$ws.Range("J80").Value = 7428.4287
$ws.Range("K80").Value = 200002600
$ws.Range("L80").Value = 7428.4287
$ws.Range("M80").Value = -200001602
$ws.Range("N80").Value = -9424.4287
$ws.Range("H83").Value = 83338750
$ws.Range("I83").Value = 200002600
$ws.Range("J83").Value = 7428.4287
$ws.Range("K83").Value = 1000013000
$ws.Range("L83").Value = 37142.14350000001
$ws.Range("M83").Value = -1000008008
$ws.Range("N83").Value = -47126.14350000001
$ws.Range("H97").Value = 486.63635
$ws.Range("I97").Value = 495.7143
$ws.Range("J97").Value = 470.75
$ws.Range("K97").Value = 495.7143
$ws.Range("L97").Value = 470.75
$ws.Range("M97").Value = 0.2857000000000198
$ws.Range("N97").Value = -1462.75

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3856.9143
$ws.Range("I40").Value = 3091.9565
$ws.Range("J40").Value = 5323.0835
$ws.Range("K40").Value = 3091.9565
$ws.Range("L40").Value = 5323.0835
$ws.Range("M40").Value = -2955.9565
$ws.Range("N40").Value = -5595.0835
$ws.Range("H43").Value = 22400
$ws.Range("I43").Value = 22400
$ws.Range("K43").Value = 22400
$ws.Range("M43").Value = -22207
$ws.Range("H76").Value = 69999.5
$ws.Range("J76").Value = 69999.5
$ws.Range("L76").Value = 69999.5
$ws.Range("N76").Value = -70675.5
$ws.Range("H79").Value = 69999.5
$ws.Range("J79").Value = 69999.5
$ws.Range("L79").Value = 69999.5
$ws.Range("N79").Value = -72339.5
$ws.Range("H106").Value = 21902.75
$ws.Range("J106").Value = 21902.75
$ws.Range("L106").Value = 21902.75
$ws.Range("N106").Value = -24426.75
$ws.Range("H122").Value = 3828.3125
$ws.Range("I122").Value = 3275.2856
$ws.Range("J122").Value = 7699.5
$ws.Range("K122").Value = 9825.856800000001
$ws.Range("L122").Value = 23098.5
$ws.Range("M122").Value = -7375.856800000001
$ws.Range("N122").Value = -27998.5
$ws.Range("H132").Value = 4065.5908
$ws.Range("I132").Value = 2592.524
$ws.Range("K132").Value = 7777.572
$ws.Range("M132").Value = -5247.572

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 5008947
$ws.Range("J33").Value = 17894
$ws.Range("L33").Value = 17894
$ws.Range("N33").Value = -18394
$ws.Range("H36").Value = 5008947
$ws.Range("J36").Value = 17894
$ws.Range("L36").Value = 17894
$ws.Range("N36").Value = -18394
$ws.Range("H122").Value = 411035.12
$ws.Range("J122").Value = 6017.364
$ws.Range("L122").Value = 18052.092
$ws.Range("N122").Value = -22952.092
$ws.Range("H136").Value = 2236.0286
$ws.Range("I136").Value = 1337.5714
$ws.Range("K136").Value = 4012.7142
$ws.Range("M136").Value = -1462.7142
$ws.Range("H141").Value = 94157.5
$ws.Range("J141").Value = 93933.336
$ws.Range("L141").Value = 93933.336
$ws.Range("N141").Value = -104293.336

Write-Host "Applied 216 value updates and 2 cell clears across 7 sheets."